$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48 (shifts existing rows 48-84 down to 49-85)
$ws.Rows("48:48").Insert()

# Populate the newly inserted row 48 with this week's data, matching the
# fixed attributes shared by every record in this subset (market, region,
# product, quality, unit, origin, classification) and the new weekly values.
$ws.Cells.Item(48, 1).Value = 1
$ws.Cells.Item(48, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(48, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(48, 4).Value = 44777
$ws.Cells.Item(48, 5).Value = 15
$ws.Cells.Item(48, 6).Value = 100112040
$ws.Cells.Item(48, 7).Value = "Cilantro"
$ws.Cells.Item(48, 8).Value = "Sin especificar"
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 300
$ws.Cells.Item(48, 11).Value = 3500
$ws.Cells.Item(48, 12).Value = 4000
$ws.Cells.Item(48, 13).Value = 3750
$ws.Cells.Item(48, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(48, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(48, 16).Value = 1875
$ws.Cells.Item(48, 17).Value = 2
$ws.Cells.Item(48, 18).Value = "Hortaliza"
